$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B3").Value = 59
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B9").Value = 1

$ws.Activate()
$ws.Range("B3").Select()
